$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New nper/pv/fv data for rows 2-32 (RRI test cases, including new failure/edge cases)
$data = @(
    @(0,   300,   400),
    @(0,   -1,    -3),
    @(1,   -1,    -3),
    @(12,  100,   10),
    @(12,  100,   -90),
    @(5,   0,     0),
    @(5,   -1,    5),
    @(5,   10,    10),
    @(2,   2,     8),
    @(2,   8,     2),
    @(2,   8,     0),
    @(2,   0,     10),
    @(12,  -5,    -6),
    @(1,   -5,    0),
    @(12,  -1,    -1),
    @(12,  300,   300),
    @(12,  300,   400),
    @(12,  300,   4000),
    @(12,  300,   40000),
    @(24,  300,   400),
    @(24,  300,   4000),
    @(24,  300,   40000),
    @(38,  300,   400),
    @(38,  300,   4000),
    @(38,  300,   40000),
    @(8,   10000, 2441880),
    @(4,   5000,  6000),
    @(4,   5000,  10000),
    @(1,   250,   275),
    @(2,   250,   500),
    @(3,   250,   880)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $triple = $data[$i]
    $ws.Cells.Item($r, 1).Value = $triple[0]
    $ws.Cells.Item($r, 2).Value = $triple[1]
    $ws.Cells.Item($r, 3).Value = $triple[2]
    $ws.Cells.Item($r, 4).Formula = "=_xlfn.RRI(A$r,B$r,C$r)"
}

# Grow the Table1 listobject to cover the new data range
$tbl = $ws.ListObjects.Item(1)
$lastRow = $startRow + $data.Count - 1
$tbl.Resize($ws.Range("A1:D$lastRow"))

# Match the author's final selection
$ws.Range("A5").Select()
